$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 'Wipro'
$ws.Range("B6").Value = ' Build the infrastructure required for optimal extraction, transformation, and loading o...'
$ws.Range("C6").Value = '3-8 Yrs'
$ws.Range("D6").Value = 'Gurgaon/Gurugram, Bangalore/Bengaluru'
$ws.Range("E6").Value = 'Not disclosed'
$ws.Range("F6").Value = '[''Snowflake'', ''AWS'', ''EC2'', ''RDS'', ''SQL'', ''IT Skills'', ''Python'', ''Cloud'']'

# Row 7
$ws.Range("A7").Value = 'Cyient'
$ws.Range("B7").Value = ' Have experience in developing automation frameworksExperience in developing automation ...'
$ws.Range("C7").Value = '3-6 Yrs'
$ws.Range("D7").Value = 'Bangalore/Bengaluru'
$ws.Range("E7").Value = 'Not disclosed'
$ws.Range("F7").Value = '[''communication'', ''troubleshooting'', ''SVN'', ''Git'', ''automation frameworks'', ''Bugzilla'', ''Python'', ''C programming'']'

# Row 8
$ws.Range("A8").Value = 'HyrEzy Talent Solutions'
$ws.Range("B8").Value = ' Graduate / Post-Graduate in Computer Science / Mathematics/Physics or allied fieldsExp ...'
$ws.Range("C8").Value = '3-6 Yrs'
$ws.Range("D8").Value = 'Mumbai'
$ws.Range("E8").Value = 'Not disclosed'
$ws.Range("F8").Value = '[''SNS'', ''ETL'', ''EC2'', ''AWS'', ''JavaScript'', ''Apache'', ''REST services'', ''Django'']'

# Row 9
$ws.Range("A9").Value = 'HyrEzy Talent Solutions'
$ws.Range("B9").Value = ' Graduate / Post-Graduate in Computer Science / Mathematics/Physics or allied fields Exp...'
$ws.Range("C9").Value = '3-6 Yrs'
$ws.Range("D9").Value = 'Mumbai'
$ws.Range("E9").Value = 'Not disclosed'
$ws.Range("F9").Value = '[''IT Skills'', ''Java'', ''Python'', ''Javascript'', ''AWS'', ''S3'', ''REST services'', ''SQS'']'

# Row 10
$ws.Range("A10").Value = 'HyrEzy Talent Solutions'
$ws.Range("B10").Value = ' Graduate / Post-Graduate in Computer Science / Mathematics/Physics or allied fieldsExp ...'
$ws.Range("C10").Value = '3-6 Yrs'
$ws.Range("D10").Value = 'Mumbai'
$ws.Range("E10").Value = 'Not disclosed'
$ws.Range("F10").Value = '[''SNS'', ''ETL'', ''EC2'', ''AWS'', ''JavaScript'', ''Apache'', ''REST services'', ''Django'']'

# Row 11
$ws.Range("A11").Value = 'Verisk Analytics India Private Limited'
$ws.Range("B11").Value = ' Minimum of bachelors degree 4 years development experienceAWS Experience is added advan...'
$ws.Range("C11").Value = '2-4 Yrs'
$ws.Range("D11").Value = 'Hyderabad/Secunderabad'
$ws.Range("E11").Value = 'Not disclosed'
$ws.Range("F11").Value = '[''written communication'', ''AJAX'', ''JavaScript'', ''XML'', ''SQL'', ''Python'', ''AWS'', ''IT Skills'']'

# Row 12
$ws.Range("A12").Value = 'Informatica'
$ws.Range("B12").Value = ' In this role, you must be able to work and adapt in a fluid, fast-paced environmentyou ...'
$ws.Range("C12").Value = '4-7 Yrs'
$ws.Range("D12").Value = 'Bangalore/Bengaluru'
$ws.Range("E12").Value = 'Not disclosed'
$ws.Range("F12").Value = '[''Statistical programming'', ''Networking'', ''Data management'', ''Access management'', ''devops'', ''Informatica'', ''SDK'', ''Vulnerability management'']'

# Row 13
$ws.Range("A13").Value = 'ANVETA CONSULTING PRIVATE LIMITED '
$ws.Range("B13").Value = ' Bachelors Degree in Computer Science or similar fieldAngular is preferred Docker / Cont...'
$ws.Range("C13").Value = '5-10 Yrs'
$ws.Range("D13").Value = 'Bangalore/Bengaluru'
$ws.Range("E13").Value = 'Not disclosed'
$ws.Range("F13").Value = '[''Fullstack Developer'', ''Django'', ''React'']'

# Row 14
$ws.Range("A14").Value = 'Societe Generale Global Solution Centre Pvt Ltd'
$ws.Range("B14").Value = ' Create and maintain a operational run book for the teamMonitor the infra capacity and m...'
$ws.Range("C14").Value = '6-7 Yrs'
$ws.Range("D14").Value = 'Bangalore/Bengaluru'
$ws.Range("E14").Value = 'Not disclosed'
$ws.Range("F14").Value = '[''IT Skills'', ''Python'', ''Cloud'', ''DevOps'', ''Jenkins'', ''AWS'', ''Azure'', ''Application Management'']'

# Row 15
$ws.Range("A15").Value = 'Onward Technologies Limited'
$ws.Range("B15").Value = ' Required key skills are mentioned below, Python & Flask FrameworkWeb application develo...'
$ws.Range("C15").Value = '6-9 Yrs'
$ws.Range("D15").Value = 'Bangalore/Bengaluru'
$ws.Range("E15").Value = '10,00,000 - 16,00,000 PA.'
$ws.Range("F15").Value = '[''GitHub'', ''SQL Server'', ''CSS'', ''Angular'', ''Flask Framework'', ''Agile methodology'', ''PostgreSQL'', ''Microservices'']'

# Row 16
$ws.Range("A16").Value = 'Future Focus Infotech Pvt. Ltd.'
$ws.Range("B16").Value = ' JOB DESCRIPTION :- Python DeveloperDjango PythonHTMLCSS ExcelHtml5'
$ws.Range("C16").Value = '3-6 Yrs'
$ws.Range("D16").Value = 'Kolkata, Hyderabad/Secunderabad, Pune, Ahmedabad, Chennai, Bangalore/Bengaluru, Delhi / NCR, Mumbai (All Areas)'
$ws.Range("E16").Value = 'Not disclosed'
$ws.Range("F16").Value = '[''CSS'', ''HTML'', ''Excel'', ''Django'', ''Html5'', ''Python'', ''IT Skills'', ''Java'']'

# Row 17
$ws.Range("A17").Value = 'Societe Generale Global Solution Centre Pvt Ltd'
$ws.Range("B17").Value = ' We are looking for a React.JS developer with working experience on Python Development t...'
$ws.Range("C17").Value = '4-9 Yrs'
$ws.Range("D17").Value = 'Bangalore/Bengaluru'
$ws.Range("E17").Value = 'Not disclosed'
$ws.Range("F17").Value = '[''Front End'', ''Html5'', ''JavaScript'', ''PHP'', ''HTML'', ''React.js'', ''Angular'', ''Python'']'

# Row 18
$ws.Range("A18").Value = 'Catalyst'
$ws.Range("B18").Value = ' Build and maintain real-time / batch data pipelines that can consolidate and clean up u...'
$ws.Range("C18").Value = '3-8 Yrs'
$ws.Range("D18").Value = 'Bangalore/Bengaluru'
$ws.Range("E18").Value = 'Not disclosed'
$ws.Range("F18").Value = '[''IT Skills'', ''Java'', ''Python'', ''Data Science'', ''Machine Learning'', ''Big Data'', ''Hive'', ''Data Pipeline'']'
